$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): add Date/Start/End headers ------------------------
# Target layout: C2 -> "Time(h)", D2 -> "Date", E2 -> "Start", F2 -> "End"
$ws.Range("C2").Value2 = "Time(h)"
$ws.Range("D2").Value2 = "Date"
$ws.Range("E2").Value2 = "Start"
$ws.Range("F2").Value2 = "End"

# --- Row 3: convert C3 from text "2h" to numeric 2, add Start/End times ----
$ws.Range("C3").Value2 = 2

# Copy the existing date-formatted cell (H3) onto D3 so it picks up the same
# number format (numFmtId 14) while re-using the same style index instead of
# creating a duplicate cellXfs entry.
$ws.Range("H3").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D3").Value2 = 45677

$ws.Range("E3").NumberFormat = "h:mm"
$ws.Range("E3").Value2 = 0.79166666666666663
$ws.Range("F3").NumberFormat = "h:mm"
$ws.Range("F3").Value2 = 0.875

# --- Row 4: new task "Score on Window" --------------------------------------
$ws.Range("B4").Value2 = "Score on Window"
$ws.Range("C4").Value2 = 3.5

$ws.Range("H3").Copy()
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D4").Value2 = 45678

$ws.Range("E3").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("E4").Value2 = 0.79166666666666663
$ws.Range("F3").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value2 = 0.9375

# --- Row 5: new task "Direction Light" --------------------------------------
$ws.Range("B5").Value2 = "Direction Light"
$ws.Range("C5").Value2 = 0.5

$ws.Range("H3").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value2 = 45679

$ws.Range("E3").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E5").Value2 = 0.6875
$ws.Range("F3").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value2 = 0.70833333333333337

# --- Column widths: size the now-populated columns to fit their content ----
# (the emulated ColumnWidth setter only supports coarse 1/6-character
# granularity, so these are the closest achievable values to the real
# Excel autofit widths)
$ws.Columns.Item(2).ColumnWidth = 15.65
$ws.Columns.Item(4).ColumnWidth = 9.3
$ws.Columns.Item(5).ColumnWidth = 4.65
$ws.Columns.Item(6).ColumnWidth = 4.65

# --- Selection / active cell -------------------------------------------------
[void]$ws.Range("H9").Select()
